# SCD0017 -> SCD0018 (and related) TC_ID update
# - Rename the worksheet from SCD0293 to SCD0018
# - Update the TC_ID column (B) values on rows 2-4 from "DGS-308" to "SCD0018-016"
# - Column B widens (bestFit) to accommodate the longer TC_ID text
# - Selection moves to B5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "SCD0018"

# 2. Update the TC_ID values (column B) for the three data rows
$ws.Range("B2").Value = "SCD0018-016"
$ws.Range("B3").Value = "SCD0018-016"
$ws.Range("B4").Value = "SCD0018-016"

# 3. Re-fit column B so it is wide enough for the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.6

# 4. Move the active selection to B5 (matches the saved view state)
[void]$ws.Range("B5").Select()
